# ------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1) The table on slide 5 (the "Type of document / Definition / Why it
#    is important" table) is switched from the deck's custom table
#    style ({078CB856-2A2D-4BE7-B0B2-4AB881732E02}) to the built-in
#    "Medium Style 2 - Accent 1" table style
#    ({759282F4-EBC5-4709-A1BC-2EDAD19B5619}).
#
# 2) The presentation's theme (ppt/theme/theme1.xml, used by the
#    slide master / all slides) is switched from the "Integral" /
#    "Red Violet" colour scheme to the default Office Theme colour
#    scheme. Table styles/colours are theme driven, so this is done
#    by writing each of the twelve theme colours through
#    Slide.ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
#    folHlink) -- the same slots PowerPoint's Design > Colors picker
#    edits.
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Re-style the table on slide 5 ------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{759282F4-EBC5-4709-A1BC-2EDAD19B5619}")

# --- 2) Swap the theme colour scheme back to the default Office ----
#        Theme colours (was the custom "Integral"/"Red Violet" set).
$themeColors = $p.Slides.Item(1).ThemeColorScheme

# index : slot      : new RGB (hex)   -> VBA RGB() packed value
#   1   : dk1        000000 -> 0
#   2   : lt1        FFFFFF -> 16777215
#   3   : dk2        44546A -> 6968388
#   4   : lt2        E7E6E6 -> 15132391
#   5   : accent1    5B9BD5 -> 13998939
#   6   : accent2    ED7D31 -> 3243501
#   7   : accent3    A5A5A5 -> 10855845
#   8   : accent4    FFC000 -> 49407
#   9   : accent5    4472C4 -> 12874308
#  10   : accent6    70AD47 -> 4697456
#  11   : hlink      0563C1 -> 12673797
#  12   : folHlink   954F72 -> 7491477
$newRgb = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $newRgb[$i - 1]
}
